$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 260.45456
$ws.Range("J2").Value = 405.66666
$ws.Range("L2").Value = 405.66666
$ws.Range("N2").Value = -631.66666
$ws.Range("H18").Value = 530.875
$ws.Range("I18").Value = 269.8
$ws.Range("J18").Value = 966
$ws.Range("K18").Value = 269.8
$ws.Range("L18").Value = 966
$ws.Range("M18").Value = 14.19999999999999
$ws.Range("N18").Value = -1534
$ws.Range("H40").Value = 7063.6924
$ws.Range("J40").Value = 7063.6924
$ws.Range("L40").Value = 7063.6924
$ws.Range("N40").Value = -7413.6924
$ws.Range("H62").Value = 11118139
$ws.Range("I62").Value = 12827854
$ws.Range("K62").Value = 12827854
$ws.Range("M62").Value = -12827230
$ws.Range("H64").Value = 55563870
$ws.Range("I64").Value = 166669120
$ws.Range("J64").Value = 11250.75
$ws.Range("K64").Value = 166669120
$ws.Range("L64").Value = 11250.75
$ws.Range("M64").Value = -166668872
$ws.Range("N64").Value = -11746.75
$ws.Range("H65").Value = 11118139
$ws.Range("I65").Value = 12827854
$ws.Range("K65").Value = 64139270
$ws.Range("M65").Value = -64136150
$ws.Range("H67").Value = 55563870
$ws.Range("I67").Value = 166669120
$ws.Range("J67").Value = 11250.75
$ws.Range("K67").Value = 166669120
$ws.Range("L67").Value = 11250.75
$ws.Range("M67").Value = -166668262
$ws.Range("N67").Value = -12966.75
$ws.Range("H107").Value = 1786.6552
$ws.Range("I107").Value = 1820.5385
$ws.Range("J107").Value = 1493
$ws.Range("K107").Value = 1820.5385
$ws.Range("L107").Value = 1493
$ws.Range("M107").Value = 99.46149999999989
$ws.Range("N107").Value = -5333
$ws.Range("H116").Value = 22224500
$ws.Range("I116").Value = 27780126
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 27780126
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -27776684
$ws.Range("N116").Value = -8884
$ws.Range("H132").Value = 250812.7
$ws.Range("I132").Value = 269774.75
$ws.Range("K132").Value = 809324.25
$ws.Range("M132").Value = -806794.25
$ws.Range("H137").Value = 5757.3335
$ws.Range("I137").Value = 13749.25
$ws.Range("J137").Value = 3473.9285
$ws.Range("K137").Value = 41247.75
$ws.Range("L137").Value = 10421.7855
$ws.Range("M137").Value = -38697.75
$ws.Range("N137").Value = -15521.7855
$ws.Range("H138").Value = 4250.25
$ws.Range("I138").Value = 3737.1428
$ws.Range("K138").Value = 11211.4284
$ws.Range("M138").Value = -6071.428400000001
$ws.Range("H141").Value = 5200
$ws.Range("I141").Value = 3150
$ws.Range("K141").Value = 9450
$ws.Range("M141").Value = -4270

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7161.2
$ws.Range("I74").Value = 8873.143
$ws.Range("J74").Value = 3166.6667
$ws.Range("K74").Value = 8873.143
$ws.Range("L74").Value = 3166.6667
$ws.Range("M74").Value = -7999.143
$ws.Range("N74").Value = -4914.6667
$ws.Range("H77").Value = 7161.2
$ws.Range("I77").Value = 8873.143
$ws.Range("J77").Value = 3166.6667
$ws.Range("K77").Value = 44365.715
$ws.Range("L77").Value = 15833.3335
$ws.Range("M77").Value = -39997.715
$ws.Range("N77").Value = -24569.3335
$ws.Range("H102").Value = 3811.5715
$ws.Range("I102").Value = 3821.9167
$ws.Range("K102").Value = 3821.9167
$ws.Range("M102").Value = -2199.9167
$ws.Range("H131").Value = 125000
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H132").Value = 778140.5
$ws.Range("I132").Value = 857775.9399999999
$ws.Range("K132").Value = 2573327.82
$ws.Range("M132").Value = -2570797.82

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8693.85
$ws.Range("J99").Value = 9073.913
$ws.Range("L99").Value = 9073.913
$ws.Range("N99").Value = -12069.913
$ws.Range("H105").Value = 2322.2273
$ws.Range("I105").Value = 2360.923
$ws.Range("K105").Value = 2360.923
$ws.Range("M105").Value = -613.9229999999998
$ws.Range("H107").Value = 4007081.2
$ws.Range("I107").Value = 5269528
$ws.Range("J107").Value = 9333
$ws.Range("K107").Value = 5269528
$ws.Range("L107").Value = 9333
$ws.Range("M107").Value = -5267608
$ws.Range("N107").Value = -13173
$ws.Range("H132").Value = 99000
$ws.Range("J132").Value = 99000
$ws.Range("L132").Value = 99000
$ws.Range("N132").Value = -109120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22731368
$ws.Range("I16").Value = 27780118
$ws.Range("J16").Value = 11999.75
$ws.Range("K16").Value = 27780118
$ws.Range("L16").Value = 11999.75
$ws.Range("M16").Value = -27779831
$ws.Range("N16").Value = -12573.75
$ws.Range("H31").Value = 11636062
$ws.Range("I31").Value = 41683916
$ws.Range("J31").Value = 4634.484
$ws.Range("K31").Value = 41683916
$ws.Range("L31").Value = 4634.484
$ws.Range("M31").Value = -41683621
$ws.Range("N31").Value = -5224.484
$ws.Range("H34").Value = 11636062
$ws.Range("I34").Value = 41683916
$ws.Range("J34").Value = 4634.484
$ws.Range("K34").Value = 41683916
$ws.Range("L34").Value = 4634.484
$ws.Range("M34").Value = -41683714
$ws.Range("N34").Value = -5038.484
$ws.Range("H58").Value = 83347690
$ws.Range("I58").Value = 111119450
$ws.Range("K58").Value = 111119450
$ws.Range("M58").Value = -111119247
$ws.Range("H113").Value = 22731368
$ws.Range("I113").Value = 27780118
$ws.Range("J113").Value = 11999.75
$ws.Range("K113").Value = 27780118
$ws.Range("L113").Value = 11999.75
$ws.Range("M113").Value = -27777948
$ws.Range("N113").Value = -16339.75
$ws.Range("H132").Value = 36089.8
$ws.Range("I132").Value = 27975
$ws.Range("J132").Value = 41499.668
$ws.Range("K132").Value = 83925
$ws.Range("L132").Value = 124499.004
$ws.Range("M132").Value = -81395
$ws.Range("N132").Value = -129559.004
$ws.Range("H134").Value = 125017096
$ws.Range("I134").Value = 200010900
$ws.Range("J134").Value = 27416.334
$ws.Range("K134").Value = 600032700
$ws.Range("L134").Value = 82249.00199999999
$ws.Range("M134").Value = -600030165
$ws.Range("N134").Value = -87319.00199999999
$ws.Range("H136").Value = 83347690
$ws.Range("I136").Value = 111119450
$ws.Range("K136").Value = 333358350
$ws.Range("M136").Value = -333355800

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 986.32355
$ws.Range("I5").Value = 623.0909
$ws.Range("J5").Value = 1652.25
$ws.Range("K5").Value = 1869.2727
$ws.Range("L5").Value = 4956.75
$ws.Range("M5").Value = -1757.2727
$ws.Range("N5").Value = -5180.75
$ws.Range("H33").Value = 698.875
$ws.Range("J33").Value = 853.5
$ws.Range("L33").Value = 5121
$ws.Range("N33").Value = -5687
$ws.Range("H94").Value = 13498.5
$ws.Range("J94").Value = 13498.5
$ws.Range("L94").Value = 40495.5
$ws.Range("N94").Value = -41847.5
$ws.Range("H135").Value = 986.32355
$ws.Range("I135").Value = 623.0909
$ws.Range("J135").Value = 1652.25
$ws.Range("K135").Value = 5607.8181
$ws.Range("L135").Value = 14870.25
$ws.Range("M135").Value = -3072.8181
$ws.Range("N135").Value = -19940.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12666.667
$ws.Range("J46").Value = 12500
$ws.Range("L46").Value = 12500
$ws.Range("N46").Value = -12812
$ws.Range("H54").Value = 3249.5
$ws.Range("J54").Value = 3249.5
$ws.Range("L54").Value = 3249.5
$ws.Range("N54").Value = -4029.5
$ws.Range("H70").Value = 5525.9414
$ws.Range("I70").Value = 5245.857
$ws.Range("K70").Value = 5245.857
$ws.Range("M70").Value = -4975.857
$ws.Range("H73").Value = 5525.9414
$ws.Range("I73").Value = 5245.857
$ws.Range("K73").Value = 5245.857
$ws.Range("M73").Value = -4309.857
$ws.Range("H93").Value = 95251
$ws.Range("J93").Value = 95251
$ws.Range("L93").Value = 95251
$ws.Range("N93").Value = -98995
$ws.Range("H102").Value = 923576.5
$ws.Range("I102").Value = 1670758.5
$ws.Range("J102").Value = 6580.364
$ws.Range("K102").Value = 1670758.5
$ws.Range("L102").Value = 6580.364
$ws.Range("M102").Value = -1669136.5
$ws.Range("N102").Value = -9824.364
$ws.Range("H113").Value = 8933.467000000001
$ws.Range("I113").Value = 5250.5
$ws.Range("J113").Value = 13142.571
$ws.Range("K113").Value = 5250.5
$ws.Range("L113").Value = 13142.571
$ws.Range("M113").Value = -3080.5
$ws.Range("N113").Value = -17482.571
$ws.Range("H126").Value = 16674446
$ws.Range("I126").Value = 26320370
$ws.Range("J126").Value = 13304
$ws.Range("K126").Value = 78961110
$ws.Range("L126").Value = 39912
$ws.Range("M126").Value = -78958640
$ws.Range("N126").Value = -44852
$ws.Range("H132").Value = 10453.125
$ws.Range("I132").Value = 10518.571
$ws.Range("J132").Value = 9995
$ws.Range("K132").Value = 31555.713
$ws.Range("L132").Value = 29985
$ws.Range("M132").Value = -29025.713
$ws.Range("N132").Value = -35045

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4654.4707
$ws.Range("J40").Value = 4390.75
$ws.Range("L40").Value = 4390.75
$ws.Range("N40").Value = -4662.75
$ws.Range("H122").Value = 1819623
$ws.Range("I122").Value = 3994971
$ws.Range("K122").Value = 11984913
$ws.Range("M122").Value = -11982463
$ws.Range("H136").Value = 12193.714
$ws.Range("I136").Value = 12892.667
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 38678.001
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -36128.001
$ws.Range("N136").Value = -29100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2349.3
$ws.Range("J81").Value = 3624.75
$ws.Range("L81").Value = 7249.5
$ws.Range("N81").Value = -9371.5
$ws.Range("H84").Value = 2349.3
$ws.Range("J84").Value = 3624.75
$ws.Range("L84").Value = 36247.5
$ws.Range("N84").Value = -46855.5
$ws.Range("H104").Value = 86818.8
$ws.Range("J104").Value = 86818.8
$ws.Range("L104").Value = 86818.8
$ws.Range("N104").Value = -93806.8
$ws.Range("H107").Value = 8334080.5
$ws.Range("I107").Value = 12500634
$ws.Range("K107").Value = 37501902
$ws.Range("M107").Value = -37499982
$ws.Range("H126").Value = 6356.769
$ws.Range("I126").Value = 3440
$ws.Range("K126").Value = 10320
$ws.Range("M126").Value = -7850
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 16142762
$ws.Range("I136").Value = 31267670
$ws.Range("K136").Value = 93803010
$ws.Range("M136").Value = -93800460
